# The document contains a single, single-column table where each row
# holds one stat value. This edit:
#   - Updates several simple value cells in place.
#   - Replaces the three "raw stats line" cells (row, tab-separated
#     numbers) near the bottom of the table with the short summary
#     values that used to live in the first three rows.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($rowIndex, $newText) {
    $cell = $t.Rows.Item($rowIndex).Cells.Item(1)
    $cell.Range.Text = $newText
}

# Simple value replacements (top of the table).
Set-CellText 1 "0M"
Set-CellText 2 "0M"
Set-CellText 3 "0M"
Set-CellText 4 "152"
Set-CellText 7 "0.01411"
Set-CellText 8 "0.00233"
Set-CellText 12 "0.76681"

# Collapse the three tab-separated raw-stat rows down to their summary
# value (these previously held "1<tab>...<tab>100.0" style runs).
Set-CellText 44 "98.42"
Set-CellText 45 "0.77"
Set-CellText 46 "48"
